$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General->Text) number format on Price (D) column cells whose
# new values would otherwise be auto-parsed as numbers by Excel, so the
# literal string representation (e.g. "40.10", "0.0000111") is preserved.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated values (prices, volume %, and the two swapped coin pairs).
$ws.Range('D2').Value = '71.069.06'
$ws.Range('E2').Value = '  +6.52%  '
$ws.Range('D3').Value = '3.679.94'
$ws.Range('E3').Value = '  +19.00%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '598.54'
$ws.Range('E5').Value = '  +3.86%  '
$ws.Range('D6').Value = '184.51'
$ws.Range('E6').Value = '  +6.73%  '
$ws.Range('D7').Value = '3.677.21'
$ws.Range('E7').Value = '  +19.01%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +4.43%  '
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  +7.79%  '
$ws.Range('D11').Value = '6.64'
$ws.Range('E11').Value = '  +4.02%  '
$ws.Range('D12').Value = '0.499'
$ws.Range('E12').Value = '  +5.60%  '
$ws.Range('D13').Value = '40.10'
$ws.Range('E13').Value = '  +12.32%  '
$ws.Range('D14').Value = '0.0000254'
$ws.Range('E14').Value = '  +6.24%  '
$ws.Range('D15').Value = '4.294.71'
$ws.Range('E15').Value = '  +19.14%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '71.165.87'
$ws.Range('E16').Value = '  +6.79%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.669.13'
$ws.Range('E17').Value = '  +18.77%  '
$ws.Range('D18').Value = '0.123'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').Value = '7.51'
$ws.Range('E19').Value = '  +7.67%  '
$ws.Range('D20').Value = '17.03'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').Value = '514.14'
$ws.Range('E21').Value = '  +6.41%  '
$ws.Range('D22').Value = '9.19'
$ws.Range('E22').Value = '  +18.36%  '
$ws.Range('D23').Value = '0.743'
$ws.Range('E23').Value = '  +7.93%  '
$ws.Range('D24').Value = '87.38'
$ws.Range('E24').Value = '  +4.87%  '
$ws.Range('D25').Value = '2.46'
$ws.Range('E25').Value = '  +11.03%  '
$ws.Range('D26').Value = '13.52'
$ws.Range('E26').Value = '  +7.07%  '
$ws.Range('D27').Value = '11.08'
$ws.Range('E27').Value = '  +10.44%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +11.90%  '
$ws.Range('D30').Value = '8.17'
$ws.Range('E30').Value = '  +2.53%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.77'
$ws.Range('E31').Value = '  +7.20%  '
$ws.Range('D32').Value = '31.62'
$ws.Range('E32').Value = '  +13.36%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').Value = '0.0000111'
$ws.Range('E33').Value = '  +18.50%  '
$ws.Range('E34').Value = '  +4.36%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '6.11'
$ws.Range('E36').Value = '  +9.53%  '
$ws.Range('E37').Value = '  +8.14%  '
$ws.Range('D38').Value = '0.345'
$ws.Range('E38').Value = '  +12.24%  '
$ws.Range('D39').Value = '2.16'
$ws.Range('E39').Value = '  +10.32%  '
$ws.Range('D40').Value = '51.01'
$ws.Range('E40').Value = '  +4.08%  '
$ws.Range('D42').Value = '45.60'
$ws.Range('E42').Value = '  -5.28%  '
$ws.Range('D43').Value = '3.158.90'
$ws.Range('E43').Value = '  +13.01%  '
$ws.Range('D44').Value = '8.83'
$ws.Range('E44').Value = '  +6.88%  '
$ws.Range('D45').Value = '412.34'
$ws.Range('E45').Value = '  +12.05%  '
$ws.Range('E46').Value = '  +6.60%  '
$ws.Range('E47').Value = '  +6.52%  '
$ws.Range('D48').Value = '28.27'
$ws.Range('E48').Value = '  +15.84%  '
$ws.Range('D49').Value = '137.71'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').Value = '2.46'
$ws.Range('E51').Value = '  +12.87%  '
